$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 111, shifting the existing row 111 (and below) down to 112.
$ws.Rows.Item(111).Insert()

# Populate the newly inserted row 111 with this week's data.
$ws.Cells.Item(111, 1).Value = 10
$ws.Cells.Item(111, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(111, 3).Value = "La Araucanía"
$ws.Cells.Item(111, 4).Value = 45239
$ws.Cells.Item(111, 5).Value = 9
$ws.Cells.Item(111, 6).Value = 100112042
$ws.Cells.Item(111, 7).Value = "Locoto"
$ws.Cells.Item(111, 8).Value = "Sin especificar"
$ws.Cells.Item(111, 9).Value = "Primera"
$ws.Cells.Item(111, 10).Value = 60
$ws.Cells.Item(111, 11).Value = 3800
$ws.Cells.Item(111, 12).Value = 3800
$ws.Cells.Item(111, 13).Value = 3800
$ws.Cells.Item(111, 14).Value = "`$/kilo"
$ws.Cells.Item(111, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(111, 16).Value = 3800
$ws.Cells.Item(111, 17).Value = 1
$ws.Cells.Item(111, 18).Value = "Hortaliza"
